# scintillator.xlsx parameter sheet update:
# "7mm scintilator, fiber slopes down to meet with center of panel"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove parameters that are no longer used ---
# Delete bottom-up so earlier row numbers stay stable while we work.

# MountingScrewDepth (row 23)
$ws.Rows("23:23").Delete()

# MPPCDepth, LongFiberAngle, ShortFiberAngle (rows 16-18)
$ws.Rows("16:18").Delete()

# MPPCClearence, MPPCSensorClearence, MPPCSensorOffsett, MPPCWidth, MPPCHeight (rows 10-14)
$ws.Rows("10:14").Delete()

# SmallTrackDiameter, LargeTrackRatio (rows 4-5)
$ws.Rows("4:5").Delete()

# --- Insert the new DeepTrackDepth parameter right after ShallowTrackDepth ---
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "DeepTrackDepth"
$ws.Range("B8").Formula = "=B2/2+(B7/2)"
$ws.Range("C8").Value = "mm"

# --- Update changed parameter values ---
# ScintillatorThickness: 10 -> 7 mm
$ws.Range("B2").Value = 7

# OpticalFiberClearence: 0.15 -> 0.1 mm
$ws.Range("B3").Value = 0.1

# MountingScrewOffsett: 20 -> 12.5 mm
$ws.Range("B13").Value = 12.5

# Restore the active cell selection left by the author
$ws.Range("B17").Select()
